# Apply the budget updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Spent" (C) value for Rent (row 2): 10000 -> 5000
$ws.Range("C2").Value = 5000

# Update the "Budgeted" (B) and "Spent" (C) values for Tithe (row 9)
$ws.Range("B9").Value = 5000
$ws.Range("C9").Value = 2000

# Totals in row 11 are formulas (SUM(B2:B9) / SUM(C2:C9)) and will
# recalculate automatically once the inputs above change.

# Move the active selection from E7 to D3, matching the saved view state.
$ws.Range("D3").Select()
